# commit: fix SF max value
#
# The existing "cop", "cop_extrapolation" and "QConMax" sheets are
# unchanged in content; a new "PEleMax" sheet is appended that derives the
# maximum electrical power (PEleMax = QConMax / COP, element-wise) from
# the data already present on "QConMax" and "cop_extrapolation", and then
# reports the overall maximum in C31.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Refresh the selection/view state on the three pre-existing sheets.
# ---------------------------------------------------------------------
$wsCop = $wb.Worksheets.Item("cop")
$wsCop.Range("D30").Select()

$wsCopExtrap = $wb.Worksheets.Item("cop_extrapolation")
$wsCopExtrap.Range("A1:J6").Select()

$wsQConMax = $wb.Worksheets.Item("QConMax")
$wsQConMax.Range("A1:J6").Select()

# ---------------------------------------------------------------------
# Append the new "PEleMax" sheet as the last tab and make it active.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $lastSheet)
$ws4.Name = "PEleMax"

# NB: "Max" (-> B31) is written before "COP" (-> A12) and "PeleMax"
# (-> A20) so the three labels land in xl/sharedStrings.xml in the order
# 0 = Max, 1 = COP, 2 = PeleMax, matching the target workbook exactly.
$ws4.Range("B31").Value = "Max"

# --- Rows 1-6: verbatim copy of the "QConMax" table --------------------
$ws4.Range("A1").Value = 0
$ws4.Range("B1").Value = 253.15
$ws4.Range("C1").Value = 258.14999999999998
$ws4.Range("D1").Value = 266.14999999999998
$ws4.Range("E1").Value = 275.14999999999998
$ws4.Range("F1").Value = 280.14999999999998
$ws4.Range("G1").Value = 283.14999999999998
$ws4.Range("H1").Value = 293.14999999999998
$ws4.Range("I1").Value = 303.14999999999998
$ws4.Range("J1").Value = 308.14999999999998
$ws4.Range("A2").Value = 308.14999999999998
$ws4.Range("B2").Value = 4490
$ws4.Range("C2").Value = 5170
$ws4.Range("D2").Value = 6470
$ws4.Range("E2").Value = 6790
$ws4.Range("F2").Value = 8000
$ws4.Range("G2").Value = 10210
$ws4.Range("H2").Value = 12330
$ws4.Range("I2").Value = 12310
$ws4.Range("J2").Value = 13090
$ws4.Range("A3").Value = 318.14999999999998
$ws4.Range("B3").Value = 4230
$ws4.Range("C3").Value = 4900
$ws4.Range("D3").Value = 6260
$ws4.Range("E3").Value = 6780
$ws4.Range("F3").Value = 8370
$ws4.Range("G3").Value = 9970
$ws4.Range("H3").Value = 11520
$ws4.Range("I3").Value = 13040
$ws4.Range("J3").Value = 12640
$ws4.Range("A4").Value = 328.15
$ws4.Range("B4").Value = 3780
$ws4.Range("C4").Value = 4710
$ws4.Range("D4").Value = 6030
$ws4.Range("E4").Value = 6830
$ws4.Range("F4").Value = 8380
$ws4.Range("G4").Value = 9940
$ws4.Range("H4").Value = 11500
$ws4.Range("I4").Value = 13070
$ws4.Range("J4").Value = 13110
$ws4.Range("A5").Value = 338.15
$ws4.Range("B5").Value = 2600
$ws4.Range("C5").Value = 3170
$ws4.Range("D5").Value = 4610
$ws4.Range("E5").Value = 6320
$ws4.Range("F5").Value = 8140
$ws4.Range("G5").Value = 9550
$ws4.Range("H5").Value = 11290
$ws4.Range("I5").Value = 12100
$ws4.Range("J5").Value = 12180
$ws4.Range("A6").Value = 343.15
$ws4.Range("B6").Value = 2300
$ws4.Range("C6").Value = 3000
$ws4.Range("D6").Value = 3830
$ws4.Range("E6").Value = 5560
$ws4.Range("F6").Value = 7600
$ws4.Range("G6").Value = 8700
$ws4.Range("H6").Value = 11290
$ws4.Range("I6").Value = 12500
$ws4.Range("J6").Value = 12590

# --- Row 12 label + rows 13-18: verbatim copy of the "cop_extrapolation"
#     (COP) table -------------------------------------------------------
$ws4.Range("A12").Value = "COP"
$ws4.Range("A13").Value = 0
$ws4.Range("B13").Value = 253.15
$ws4.Range("C13").Value = 258.14999999999998
$ws4.Range("D13").Value = 266.14999999999998
$ws4.Range("E13").Value = 275.14999999999998
$ws4.Range("F13").Value = 280.14999999999998
$ws4.Range("G13").Value = 283.14999999999998
$ws4.Range("H13").Value = 293.14999999999998
$ws4.Range("I13").Value = 303.14999999999998
$ws4.Range("J13").Value = 308.14999999999998
$ws4.Range("A14").Value = 308.14999999999998
$ws4.Range("B14").Value = 2.1099624060150375
$ws4.Range("C14").Value = 2.3100983020554064
$ws4.Range("D14").Value = 2.7003338898163607
$ws4.Range("E14").Value = 3.7002724795640325
$ws4.Range("F14").Value = 4.8989589712186161
$ws4.Range("G14").Value = 5.0897308075772685
$ws4.Range("H14").Value = 7.1686046511627906
$ws4.Range("I14").Value = 8.4027303754266214
$ws4.Range("J14").Value = 8.4017971758664949
$ws4.Range("A15").Value = 318.14999999999998
$ws4.Range("B15").Value = 1.8399304045237059
$ws4.Range("C15").Value = 2.0399666944213157
$ws4.Range("D15").Value = 2.3703142748958728
$ws4.Range("E15").Value = 3.1201104463874829
$ws4.Range("F15").Value = 3.7299465240641712
$ws4.Range("G15").Value = 4.0495532087733546
$ws4.Range("H15").Value = 5.459715639810427
$ws4.Range("I15").Value = 7.520184544405998
$ws4.Range("J15").Value = 7.9797979797979801
$ws4.Range("A16").Value = 328.15
$ws4.Range("B16").Value = 1.6399132321041214
$ws4.Range("C16").Value = 1.7901938426453821
$ws4.Range("D16").Value = 2.080027595722663
$ws4.Range("E16").Value = 2.7396710790212597
$ws4.Range("F16").Value = 3.1397527163731733
$ws4.Range("G16").Value = 3.4099485420240137
$ws4.Range("H16").Value = 4.4992175273865413
$ws4.Range("I16").Value = 6.0397412199630311
$ws4.Range("J16").Value = 6.2998558385391634
$ws4.Range("A17").Value = 338.15
$ws4.Range("B17").Value = 1
$ws4.Range("C17").Value = 1.5600393700787401
$ws4.Range("D17").Value = 1.8603712671509283
$ws4.Range("E17").Value = 2.2898550724637681
$ws4.Range("F17").Value = 2.630048465266559
$ws4.Range("G17").Value = 2.8296296296296295
$ws4.Range("H17").Value = 3.6596434359805512
$ws4.Range("I17").Value = 4.7996826656088851
$ws4.Range("J17").Value = 4.8506571087216246
$ws4.Range("A18").Value = 343.15
$ws4.Range("B18").Value = 1
$ws4.Range("C18").Value = 1
$ws4.Range("D18").Value = 1.7401181281235802
$ws4.Range("E18").Value = 2.0600222304557243
$ws4.Range("F18").Value = 2.3802067021609772
$ws4.Range("G18").Value = 2.5603296056503826
$ws4.Range("H18").Value = 3.3903903903903903
$ws4.Range("I18").Value = 4.4294826364280651
$ws4.Range("J18").Value = 4.469293574724885

# --- Row 20 label + row 21 header ---------------------------------------
$ws4.Range("A20").Value = "PeleMax"
$ws4.Range("A21").Value = 0
$ws4.Range("B21").Value = 253.15
$ws4.Range("C21").Value = 258.14999999999998
$ws4.Range("D21").Value = 266.14999999999998
$ws4.Range("E21").Value = 275.14999999999998
$ws4.Range("F21").Value = 280.14999999999998
$ws4.Range("G21").Value = 283.14999999999998
$ws4.Range("H21").Value = 293.14999999999998
$ws4.Range("I21").Value = 303.14999999999998
$ws4.Range("J21").Value = 308.14999999999998

# --- Rows 22-26: row labels (temperature) -------------------------------
$ws4.Range("A22").Value = 308.14999999999998
$ws4.Range("A23").Value = 318.14999999999998
$ws4.Range("A24").Value = 328.15
$ws4.Range("A25").Value = 338.15
$ws4.Range("A26").Value = 343.15

# --- Rows 22-26: PEleMax = QConMax / COP, element-wise -------------------
$ws4.Range("B22").Formula = "=B2/B14"
$ws4.Range("C22").Formula = "=C2/C14"
$ws4.Range("D22").Formula = "=D2/D14"
$ws4.Range("E22").Formula = "=E2/E14"
$ws4.Range("F22").Formula = "=F2/F14"
$ws4.Range("G22").Formula = "=G2/G14"
$ws4.Range("H22").Formula = "=H2/H14"
$ws4.Range("I22").Formula = "=I2/I14"
$ws4.Range("J22").Formula = "=J2/J14"
$ws4.Range("B23").Formula = "=B3/B15"
$ws4.Range("C23").Formula = "=C3/C15"
$ws4.Range("D23").Formula = "=D3/D15"
$ws4.Range("E23").Formula = "=E3/E15"
$ws4.Range("F23").Formula = "=F3/F15"
$ws4.Range("G23").Formula = "=G3/G15"
$ws4.Range("H23").Formula = "=H3/H15"
$ws4.Range("I23").Formula = "=I3/I15"
$ws4.Range("J23").Formula = "=J3/J15"
$ws4.Range("B24").Formula = "=B4/B16"
$ws4.Range("C24").Formula = "=C4/C16"
$ws4.Range("D24").Formula = "=D4/D16"
$ws4.Range("E24").Formula = "=E4/E16"
$ws4.Range("F24").Formula = "=F4/F16"
$ws4.Range("G24").Formula = "=G4/G16"
$ws4.Range("H24").Formula = "=H4/H16"
$ws4.Range("I24").Formula = "=I4/I16"
$ws4.Range("J24").Formula = "=J4/J16"
$ws4.Range("B25").Formula = "=B5/B17"
$ws4.Range("C25").Formula = "=C5/C17"
$ws4.Range("D25").Formula = "=D5/D17"
$ws4.Range("E25").Formula = "=E5/E17"
$ws4.Range("F25").Formula = "=F5/F17"
$ws4.Range("G25").Formula = "=G5/G17"
$ws4.Range("H25").Formula = "=H5/H17"
$ws4.Range("I25").Formula = "=I5/I17"
$ws4.Range("J25").Formula = "=J5/J17"
$ws4.Range("B26").Formula = "=B6/B18"
$ws4.Range("C26").Formula = "=C6/C18"
$ws4.Range("D26").Formula = "=D6/D18"
$ws4.Range("E26").Formula = "=E6/E18"
$ws4.Range("F26").Formula = "=F6/F18"
$ws4.Range("G26").Formula = "=G6/G18"
$ws4.Range("H26").Formula = "=H6/H18"
$ws4.Range("I26").Formula = "=I6/I18"
$ws4.Range("J26").Formula = "=J6/J18"

# --- Row 31: overall maximum (B31 label was already written above) ------
$ws4.Range("C31").Formula = "=MAX(B22:J27)"

# --- Final view state on the new sheet ------------------------------------
$ws4.Range("A21").Select()
$ws4.Activate()
